$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells D1, E1 - set values then copy the header style from
# the existing header row (A1:B1 uses style index 1: bold, bordered,
# centered) onto the two new header cells.
$ws.Range("D1").Value = "mean_E"
$ws.Range("E1").Value = "mean_UPL_1"
$ws.Range("A1:B1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)

# New data columns D (mean_E) and E (mean_UPL_1) for rows 2-7.
$ws.Range("D2").Value = 59652.66002461074
$ws.Range("E2").Value = 163.2461647478644

$ws.Range("D3").Value = 59855.77924476722
$ws.Range("E3").Value = 145.2313242222233

$ws.Range("D4").Value = 49741.58456709496
$ws.Range("E4").Value = 130.8866971883324

$ws.Range("D5").Value = 48391.35342360516
$ws.Range("E5").Value = 134.8347027374249

$ws.Range("D6").Value = 40860.56796752069
$ws.Range("E6").Value = 86.38721691376323

$ws.Range("D7").Value = 43041.85017422398
$ws.Range("E7").Value = 56.33522942889576
